# Update column G ("K") values on Sheet1, rows 2-85,
# regenerated using strikeout counts (K) instead of the old "Strike#" values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @(1,1,2,1,2,1,0,3,1,1,1,1,2,2,1,0,2,1,1,0,2,1,1,0,2,2,0,2,0,1,3,2,1,3,2,1,1,1,1,1,1,0,0,1,3,2,2,2,2,2,0,1,1,1,3,2,0,1,2,2,2,2,2,1,1,2,3,1,2,2,0,1,3,3,2,0,1,2,0,1,2,2,1,1)

$rowCount = $kValues.Length
$arr = New-Object 'object[,]' $rowCount,1
for ($i = 0; $i -lt $rowCount; $i++) {
    $arr[$i, 0] = $kValues[$i]
}

$ws.Range("G2:G85").Value = $arr
